$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.903.28'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.898.92'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7979'
$ws.Range('E5').Value = '  -5.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.43'
$ws.Range('E6').Value = '  +0.89%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9997'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3177'
$ws.Range('E8').Value = '  -3.90%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.57'
$ws.Range('E9').Value = '  -4.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07206'
$ws.Range('E10').Value = '  +1.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08114'
$ws.Range('E11').Value = '  +0.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7738'
$ws.Range('E12').Value = '  +1.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.545'
$ws.Range('E13').Value = '  +5.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.917.98'
$ws.Range('E14').Value = '  +1.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.93'
$ws.Range('E15').Value = '  +0.52%  '
$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.154'
$ws.Range('E16').Value = '  +4.23%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.909.35'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.98'
$ws.Range('E18').Value = '  -1.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '246.19'
$ws.Range('E19').Value = '  +0.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007770'
$ws.Range('E20').Value = '  -0.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.270'
$ws.Range('E21').Value = '  +18.05%  '
$ws.Range('E22').Value = '  +0.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.148.97'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  +0.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1693'
$ws.Range('E25').Value = '  -2.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.481'
$ws.Range('E26').Value = '  +2.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '164.74'
$ws.Range('E27').Value = '  -0.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.77'
$ws.Range('E28').Value = '  -0.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.081'
$ws.Range('E29').Value = '  -1.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.415'
$ws.Range('E30').Value = '  +4.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.548'
$ws.Range('E31').Value = '  +2.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.514'
$ws.Range('E32').Value = '  +4.71%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05653'
$ws.Range('E33').Value = '  -4.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.090'
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.291'
$ws.Range('E35').Value = '  +0.98%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7449'
$ws.Range('E36').Value = '  +1.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.005'
$ws.Range('E37').Value = '  +0.74%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.630'
$ws.Range('E38').Value = '  -2.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01938'
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.781'
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.180.58'
$ws.Range('E41').Value = '  +16.72%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '75.08'
$ws.Range('E42').Value = '  +3.30%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4438'
$ws.Range('E43').Value = '  -0.44%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.913'
$ws.Range('E44').Value = '  +0.32%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8556'
$ws.Range('E45').Value = '  +1.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '104.44'
$ws.Range('E46').Value = '  +2.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9991'
$ws.Range('E47').Value = '  +0.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.891'
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.06'
$ws.Range('E49').Value = '  +2.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.492'
$ws.Range('E50').Value = '  -1.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.052.31'
$ws.Range('E51').Value = '  +0.29%  '
